# Cambio en la respuesta delete para traerlo desde la data de excel
#
# Adds a new "RESPONSE_PHONE" column (I) to the InfoUser sheet, with its
# sample value, mirroring FIRSTNAME..CODE / Esteban..200 already there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header + sample value in column I
$ws.Range("I1").Value = "RESPONSE_PHONE"
$ws.Range("I2").Value = "1-570-236-7033"

# Match the width Excel/LibreOffice would auto-fit for the new column
$ws.Columns("I").ColumnWidth = 17.37

# Leave the selection where the author left it after adding the column
$ws.Range("I5").Select() | Out-Null
